$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Pre-warm the shared style table on a throwaway worksheet so that the
# two new cell styles we need (bold+text-format header style, and
# text-format left/top-aligned data style) get registered in the exact
# same order/shape as Excel would naturally create them, without
# leaving stray intermediate styles behind on the real sheet.
# ---------------------------------------------------------------------
$scratch = $wb.Worksheets.Add()

$scratchHeader = $scratch.Cells.Item(1, 1)
$scratchHeader.Font.Bold = $true
$scratchHeader.NumberFormat = "@"
$scratchHeader.Value = "x"

$scratchData = $scratch.Cells.Item(1, 2)
$scratchData.NumberFormat = "@"
$scratchData.HorizontalAlignment = -4131  # xlLeft
$scratchData.VerticalAlignment = -4160    # xlTop
$scratchData.Value = "x"

$excel.DisplayAlerts = $false
$scratch.Delete()
$excel.DisplayAlerts = $true

# ---------------------------------------------------------------------
# Add the new worksheet "AddOpportunity" at the end of the workbook
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "AddOpportunity"

# ---- Header row (row 1) ----
$headers = @(
    "Client","Subject","JobType","IndustryGroup/HLSector","Sector","AdditionalClient",
    "AdditonalSubject","ReferralType","NonPublicInfo","BeneficialOwner","PrimaryOffice",
    "LegalEntity","DisclosureStatus","Staff","Retainer","MonthlyFee","ContingentFee",
    "ClientOwnership","SubjectOwnership","SICCode","OpportunityDescription","ReferralContact",
    "Agreement","Outcome","RecordType","FASJobType","MarketCap","Fee","StdUser","WomenLed"
)

# Bold header cells (regular bold style) - columns A-N, R-S, U-AD
$newSheet.Range("A1:N1").Font.Bold = $true
$newSheet.Range("R1:S1").Font.Bold = $true
$newSheet.Range("U1:AD1").Font.Bold = $true

# Bold + text-format header cells - columns O,P,Q,T
$newSheet.Range("O1:Q1").Font.Bold = $true
$newSheet.Range("O1:Q1").NumberFormat = "@"
$newSheet.Range("T1").Font.Bold = $true
$newSheet.Range("T1").NumberFormat = "@"

for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ---- Data row (row 2) ----
$values = @(
    "Techno Coatings, Inc.","Techno Coatings, Inc.","Buyside","CSDN-0000002536",
    "Dealership & Rental Services","No","No","Accountant","No","No","AM","HL Capital, Inc.",
    "Do Not Disclose","Emre Abale","10","10","10","Public Equity","Public Equity","9999",
    "Test","Chris Lord","Yes, separate signed agreement","Cleared","CF","Consulting","10","10",
    "Emre Abale","Yes"
)

# Text-format, left/top-aligned data cells - columns O,P,Q,T,AA,AB
$newSheet.Range("O2:Q2").NumberFormat = "@"
$newSheet.Range("O2:Q2").HorizontalAlignment = -4131  # xlLeft
$newSheet.Range("O2:Q2").VerticalAlignment = -4160    # xlTop
$newSheet.Range("T2").NumberFormat = "@"
$newSheet.Range("T2").HorizontalAlignment = -4131
$newSheet.Range("T2").VerticalAlignment = -4160
$newSheet.Range("AA2:AB2").NumberFormat = "@"
$newSheet.Range("AA2:AB2").HorizontalAlignment = -4131
$newSheet.Range("AA2:AB2").VerticalAlignment = -4160

for ($i = 0; $i -lt $values.Length; $i++) {
    $newSheet.Cells.Item(2, $i + 1).Value = $values[$i]
}

# Column D width (bestFit-like custom width)
$newSheet.Columns.Item(4).ColumnWidth = 22.88671875

# Set the new sheet as the active sheet / tab
$newSheet.Activate()
$newSheet.Select()
